$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.673.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.76%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.893.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.23%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -1.22%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'312.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.29%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  -1.12%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4875"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.97%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3794"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.76%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07335"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.62%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.9155"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.59%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'20.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.34%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.07693"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.61%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.924.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.56%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.477"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.36%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'6.617"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.14%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'91.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.06%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "'  -1.19%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000008797"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.82%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  -1.01%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'27.708.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.66%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  -2.43%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.128"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.03%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'2.131.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.87%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  -0.96%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'1.904"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.60%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'153.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.15%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'18.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.06%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'2.142"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.22%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'115.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.24%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'4.908"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.55%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.08919"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.00%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'3.186"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.50%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  -0.01%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.7671"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.23%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'4.638"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.32%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.02038"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.48%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.525"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -6.99%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -3.62%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = "'  -1.81%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.5474"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.07%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'2.984"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.40%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'6.915"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.02%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'8.484"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.21%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.1517"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -1.12%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'111.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.55%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'10.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.25%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.4798"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.99%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -1.17%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  -2.23%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'67.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.12%  "
$ws.Range("E50").Style = "Normal"
